$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.816.06'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '1.757.54'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").Value = '''237.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").Value = '''0.5077'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.74%  '
$ws.Range("D8").Value = '''41.19'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '''0.2655'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.61%  '
$ws.Range("D10").Value = '''0.06211'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.71%  '
$ws.Range("D11").Value = '1.754.67'
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").Value = '''0.06934'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.02%  '
$ws.Range("D13").Value = '''15.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.82%  '
$ws.Range("D14").Value = '''0.6062'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.78%  '
$ws.Range("D15").Value = '''4.458'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.67%  '
$ws.Range("D16").Value = '''77.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.24%  '
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("D19").Value = '25.859.34'
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").Value = '''0.000006853'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +10.24%  '
$ws.Range("E21").Value = '  +6.06%  '
$ws.Range("D22").Value = '1.976.44'
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").Value = '''4.071'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.73%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '''5.184'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.98%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '''8.144'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.11%  '
$ws.Range("D26").Value = '''137.85'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.14%  '
$ws.Range("D27").Value = '''1.458'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.60%  '
$ws.Range("D28").Value = '''1.822'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.70%  '
$ws.Range("E29").Value = '  +5.93%  '
$ws.Range("D30").Value = '''102.72'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.25%  '
$ws.Range("D31").Value = '''0.08238'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("D32").Value = '''3.704'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.68%  '
$ws.Range("D33").Value = '''3.406'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.76%  '
$ws.Range("D34").Value = '''0.04377'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.22%  '
$ws.Range("D35").Value = '''1.0000'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.34%  '
$ws.Range("E36").Value = '  +1.48%  '
$ws.Range("D37").Value = '''1.003'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("D38").Value = '''0.6008'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.81%  '
$ws.Range("D39").Value = '''2.730'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.33%  '
$ws.Range("D40").Value = '''0.01552'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.11%  '
$ws.Range("E41").Value = '  -6.81%  '
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("D43").Value = '''103.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("D44").Value = '''0.3842'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.06%  '
$ws.Range("D45").Value = '''0.7450'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.15%  '
$ws.Range("D46").Value = '''4.875'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.02%  '
$ws.Range("D47").Value = '''0.05494'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.28%  '
$ws.Range("D48").Value = '''0.1082'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.71%  '
$ws.Range("D49").Value = '''5.968'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.62%  '
$ws.Range("D50").Value = '''30.11'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.48%  '
$ws.Range("D51").Value = '''52.10'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.07%  '
